$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date stamp to cell C1 (serial 44307 = 2021-04-21), formatted as a date.
# Setting the NumberFormat before the Value avoids Excel auto-assigning a
# separate custom date format when the DateTime value is first written.
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
